$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.866.52'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.62%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.474.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.93%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '414.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.54%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.22'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.92%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.10%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.728'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.20%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +9.93%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.59'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.17%  '

# Row 12
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000230'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.38%  '

# Row 13
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.71'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +6.14%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.019.85'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.63%  '

# Row 15
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.17%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.80%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.452.97'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.56%  '

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.26%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.37%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '62.741.66'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.38%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '463.72'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.77%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '90.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.45%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.29'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.69%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.28'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.64%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.77'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +16.02%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.31'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.04%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.44'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.55%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.79'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.34%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.57'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.20%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.09'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.57%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.86%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.01%  '

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.16%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.99'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.84%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.03%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.45'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +8.58%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0492'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.65%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.10'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.45%  '

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.02%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '149.32'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.03%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.72'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.78%  '

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.322'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.48%  '

# Row 43
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.134'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.15%  '

# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.18%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.51%  '

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.73%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₃0572'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +36.69%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +11.42%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.39'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.00%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.28'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.11%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.140'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.42%  '
